$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1178.6666
$ws.Range("I74").Value = 1223.12
$ws.Range("J74").Value = 1039.75
$ws.Range("K74").Value = 1223.12
$ws.Range("L74").Value = 1039.75
$ws.Range("M74").Value = -349.1199999999999
$ws.Range("N74").Value = -2787.75

$ws.Range("H77").Value = 1178.6666
$ws.Range("I77").Value = 1223.12
$ws.Range("J77").Value = 1039.75
$ws.Range("K77").Value = 6115.599999999999
$ws.Range("L77").Value = 5198.75
$ws.Range("M77").Value = -1747.599999999999
$ws.Range("N77").Value = -13934.75

$ws.Range("H132").Value = 1618.8868
$ws.Range("I132").Value = 1487.1538
$ws.Range("J132").Value = 1985.8572
$ws.Range("K132").Value = 4461.4614
$ws.Range("L132").Value = 5957.571599999999
$ws.Range("M132").Value = -1931.4614
$ws.Range("N132").Value = -11017.5716

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 28673.365
$ws.Range("I134").Value = 1629.6415
$ws.Range("J134").Value = 96926.57000000001
$ws.Range("K134").Value = 4888.9245
$ws.Range("L134").Value = 290779.71
$ws.Range("M134").Value = -2353.9245
$ws.Range("N134").Value = -295849.71

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1838873.2
$ws.Range("I31").Value = 2370617.8
$ws.Range("J31").Value = 1937.6364
$ws.Range("K31").Value = 2370617.8
$ws.Range("L31").Value = 1937.6364
$ws.Range("M31").Value = -2370322.8
$ws.Range("N31").Value = -2527.6364

$ws.Range("H34").Value = 1838873.2
$ws.Range("I34").Value = 2370617.8
$ws.Range("J34").Value = 1937.6364
$ws.Range("K34").Value = 2370617.8
$ws.Range("L34").Value = 1937.6364
$ws.Range("M34").Value = -2370415.8
$ws.Range("N34").Value = -2341.6364

$ws.Range("H132").Value = 1562.6104
$ws.Range("I132").Value = 934.7659
$ws.Range("J132").Value = 2546.2334
$ws.Range("K132").Value = 2804.2977
$ws.Range("L132").Value = 7638.7002
$ws.Range("M132").Value = -274.2977000000001
$ws.Range("N132").Value = -12698.7002

$ws.Range("H134").Value = 5953478.5
$ws.Range("I134").Value = 1005.28986
$ws.Range("J134").Value = 33334856
$ws.Range("K134").Value = 3015.86958
$ws.Range("L134").Value = 100004568
$ws.Range("M134").Value = -480.86958
$ws.Range("N134").Value = -100009638

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 2670.6667
$ws.Range("I62").Value = 2006
$ws.Range("J62").Value = 4000
$ws.Range("K62").Value = 6018
$ws.Range("L62").Value = 12000
$ws.Range("M62").Value = -5332
$ws.Range("N62").Value = -13372

$ws.Range("H65").Value = 2670.6667
$ws.Range("I65").Value = 2006
$ws.Range("J65").Value = 4000
$ws.Range("K65").Value = 18054
$ws.Range("L65").Value = 36000
$ws.Range("M65").Value = -14622
$ws.Range("N65").Value = -42864

$ws.Range("H70").Value = 3788.875
$ws.Range("I70").Value = 3951.8333
$ws.Range("J70").Value = 3300
$ws.Range("K70").Value = 11855.4999
$ws.Range("L70").Value = 9900
$ws.Range("M70").Value = -11540.4999
$ws.Range("N70").Value = -10530

$ws.Range("H73").Value = 3788.875
$ws.Range("I73").Value = 3951.8333
$ws.Range("J73").Value = 3300
$ws.Range("K73").Value = 11855.4999
$ws.Range("L73").Value = 9900
$ws.Range("M73").Value = -10763.4999
$ws.Range("N73").Value = -12084

$ws.Range("H74").Value = 4908.3335
$ws.Range("J74").Value = 4908.3335
$ws.Range("L74").Value = 14725.0005
$ws.Range("N74").Value = -16847.0005

$ws.Range("H75").Value = 900
$ws.Range("J75").Value = 900
$ws.Range("L75").Value = 2700
$ws.Range("N75").Value = -4696

$ws.Range("H76").Value = 4001.8572
$ws.Range("I76").Value = 1506.5
$ws.Range("J76").Value = 5000
$ws.Range("K76").Value = 4519.5
$ws.Range("L76").Value = 15000
$ws.Range("M76").Value = -4136.5
$ws.Range("N76").Value = -15766

$ws.Range("H77").Value = 4908.3335
$ws.Range("J77").Value = 4908.3335
$ws.Range("L77").Value = 44175.0015
$ws.Range("N77").Value = -54783.0015

$ws.Range("H78").Value = 900
$ws.Range("J78").Value = 900
$ws.Range("L78").Value = 8100
$ws.Range("N78").Value = -18084

$ws.Range("H79").Value = 4001.8572
$ws.Range("I79").Value = 1506.5
$ws.Range("J79").Value = 5000
$ws.Range("K79").Value = 4519.5
$ws.Range("L79").Value = 15000
$ws.Range("M79").Value = -3193.5
$ws.Range("N79").Value = -17652

$ws.Range("H80").Value = 3900.1765
$ws.Range("I80").Value = 2833.3333
$ws.Range("J80").Value = 4128.7856
$ws.Range("K80").Value = 8499.999899999999
$ws.Range("L80").Value = 12386.3568
$ws.Range("M80").Value = -7563.999899999999
$ws.Range("N80").Value = -14258.3568

$ws.Range("H81").Value = 2200
$ws.Range("I81").Value = 1750
$ws.Range("J81").Value = 2500
$ws.Range("K81").Value = 5250
$ws.Range("L81").Value = 7500
$ws.Range("M81").Value = -4127
$ws.Range("N81").Value = -9746

$ws.Range("H82").Value = 12128.5
$ws.Range("I82").Value = 500
$ws.Range("J82").Value = 13023
$ws.Range("K82").Value = 1500
$ws.Range("L82").Value = 39069
$ws.Range("M82").Value = -1094
$ws.Range("N82").Value = -39881

$ws.Range("H83").Value = 3900.1765
$ws.Range("I83").Value = 2833.3333
$ws.Range("J83").Value = 4128.7856
$ws.Range("K83").Value = 25499.9997
$ws.Range("L83").Value = 37159.0704
$ws.Range("M83").Value = -20819.9997
$ws.Range("N83").Value = -46519.0704

$ws.Range("H84").Value = 2200
$ws.Range("I84").Value = 1750
$ws.Range("J84").Value = 2500
$ws.Range("K84").Value = 15750
$ws.Range("L84").Value = 22500
$ws.Range("M84").Value = -10134
$ws.Range("N84").Value = -33732

$ws.Range("H85").Value = 12128.5
$ws.Range("I85").Value = 500
$ws.Range("J85").Value = 13023
$ws.Range("K85").Value = 1500
$ws.Range("L85").Value = 39069
$ws.Range("M85").Value = -96
$ws.Range("N85").Value = -41877

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1825.3636
$ws.Range("I132").Value = 1571.0476
$ws.Range("J132").Value = 2647
$ws.Range("K132").Value = 4713.142800000001
$ws.Range("L132").Value = 7941
$ws.Range("M132").Value = -2183.142800000001
$ws.Range("N132").Value = -13001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2140.7878
$ws.Range("I132").Value = 1780.931
$ws.Range("J132").Value = 4749.75
$ws.Range("K132").Value = 5342.793
$ws.Range("L132").Value = 14249.25
$ws.Range("M132").Value = -2812.793
$ws.Range("N132").Value = -19309.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1122.3
$ws.Range("I132").Value = 1038.4286
$ws.Range("J132").Value = 1215
$ws.Range("K132").Value = 3115.2858
$ws.Range("L132").Value = 3645
$ws.Range("M132").Value = -585.2857999999997
$ws.Range("N132").Value = -8705

$ws.Range("H136").Value = 418.63333
$ws.Range("I136").Value = 373.5
$ws.Range("J136").Value = 599.1667
$ws.Range("K136").Value = 1120.5
$ws.Range("L136").Value = 1797.5001
$ws.Range("M136").Value = 1429.5
$ws.Range("N136").Value = -6897.5001
